$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.551.37"
$ws.Range("E2").Value = "  -0.45%  "

# Row 3
$ws.Range("D3").Value = "2.659.23"
$ws.Range("E3").Value = "  -0.83%  "

# Row 4
$ws.Range("E4").Value = "  -0.02%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.74"
$ws.Range("E5").Value = "  -0.75%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.88"
$ws.Range("E6").Value = "  -0.77%  "

# Row 7
$ws.Range("E7").Value = "  -0.01%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.627"
$ws.Range("E8").Value = "  +6.79%  "

# Row 9
$ws.Range("E9").Value = "  +3.18%  "

# Row 10
$ws.Range("E10").Value = "  -0.58%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.80"
$ws.Range("E11").Value = "  -1.95%  "

# Row 12
$ws.Range("E12").Value = "  +0.33%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "28.77"
$ws.Range("E13").Value = "  -2.69%  "

# Row 14
$ws.Range("E14").Value = "  -2.48%  "

# Row 15
$ws.Range("E15").Value = "  -0.88%  "

# Row 16
$ws.Range("D16").Value = "65.400.86"
$ws.Range("E16").Value = "  -0.39%  "

# Row 17
$ws.Range("D17").Value = "2.649.93"
$ws.Range("E17").Value = "  -1.65%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "12.66"
$ws.Range("E18").Value = "  +0.22%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.77"
$ws.Range("E19").Value = "  -1.25%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.48"
$ws.Range("E20").Value = "  -1.97%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "348.87"
$ws.Range("E21").Value = "  -0.86%  "

# Row 22
$ws.Range("E22").Value = "  -0.04%  "

# Row 23
$ws.Range("E23").Value = "  -2.04%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000113"
$ws.Range("E24").Value = "  +1.61%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.64"
$ws.Range("E25").Value = "  -1.80%  "

# Row 26
$ws.Range("E26").Value = "  +2.75%  "

# Row 27
$ws.Range("E27").Value = "  -2.56%  "

# Row 28
$ws.Range("E28").Value = "  -3.11%  "

# Row 29
$ws.Range("B29").Value = "Aptos"
$ws.Range("C29").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.97"
$ws.Range("E29").Value = "  -2.62%  "

# Row 30
$ws.Range("B30").Value = "Binance-PegBSC-USD"
$ws.Range("C30").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.00"
$ws.Range("E30").Value = "  +0.02%  "

# Row 31
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "534.81"
$ws.Range("E31").Value = "  +0.45%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.11"
$ws.Range("E32").Value = "  -3.08%  "

# Row 33
$ws.Range("E33").Value = "  -0.62%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.42"
$ws.Range("E34").Value = "  -2.48%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.44"
$ws.Range("E35").Value = "  +0.26%  "

# Row 36
$ws.Range("E36").Value = "  -1.73%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "20.39"
$ws.Range("E37").Value = "  -0.72%  "

# Row 38
$ws.Range("E38").Value = "  -0.09%  "

# Row 39
$ws.Range("B39").Value = "Stacks"
$ws.Range("C39").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.92"
$ws.Range("E39").Value = "  -2.61%  "

# Row 40
$ws.Range("B40").Value = "Monero"
$ws.Range("C40").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "154.93"
$ws.Range("E40").Value = "  -3.65%  "

# Row 41
$ws.Range("E41").Value = "  +0.02%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "161.35"
$ws.Range("E42").Value = "  -2.80%  "

# Row 43
$ws.Range("E43").Value = "  -0.69%  "

# Row 44
$ws.Range("E44").Value = "  +3.07%  "

# Row 45
$ws.Range("E45").Value = "  -2.68%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "22.52"
$ws.Range("E46").Value = "  -3.03%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.637"
$ws.Range("E47").Value = "  -2.14%  "

# Row 48
$ws.Range("E48").Value = "  -2.97%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0993"
$ws.Range("E49").Value = "  +0.73%  "

# Row 50
$ws.Range("D50").Value = "0.0₆0249"
$ws.Range("E50").Value = "  +5.46%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.69"
$ws.Range("E51").Value = "  -3.45%  "

# Reset number formatting on cells forced to text so style stays default
foreach ($addr in @("D5","D6","D8","D11","D13","D18","D19","D20","D21","D24","D25","D29","D30","D31","D32","D34","D35","D37","D39","D40","D42","D46","D47","D49","D51")) {
    $ws.Range($addr).ClearFormats()
}
